$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(
        "I have 20000 modifiers added ty log, why I can't I add anymore?",
        "llama3.2:latest",
        "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file. The maximum number of layouts allowed is 100 tables in a log."
    ),
    @(
        "I have 20000 modifiers added ty log, why I can't I add anymore?",
        "llama3.2:latest",
        "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file. The maximum number of layouts allowed is 20000."
    ),
    @(
        "Why can I not add another layout to my log?",
        "llama3.2:latest",
        "You cannot add another layout to your log because you have reached the limit of 19 layouts per ODF file, as specified by the GEO application."
    )
)

$startRow = 279
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
